$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Rename the existing sheet and add the three new ones, in the right order
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"

$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"

# ---------------------------------------------------------------------------
# Sheet 1 - "Sales vs PO": insert a new "Order Week" column (C) between the
# existing "y" (B) and "PO_Requested_Qty" columns, and refresh every row's
# dates/values to the corrected forecast output.
# ---------------------------------------------------------------------------
$ws1.Range("C1").Value = "Order Week"
$ws1.Range("D1").Value = "PO_Requested_Qty"

$sheet1Data = @(
  @(44920,35,44914,0),
  @(44920,0,44914,0),
  @(44927,23,44921,0),
  @(44927,0,44921,0),
  @(44983,196,44977,0),
  @(44990,86,44984,0),
  @(44997,8,44991,0),
  @(45004,5,44998,0),
  @(45011,9,45005,0),
  @(45018,1,45012,0),
  @(45025,1,45019,0),
  @(45032,2,45026,0),
  @(45039,3,45033,0),
  @(45151,7,45145,0),
  @(45158,1,45152,0),
  @(45165,2,45159,0),
  @(45214,0,45208,0),
  @(45256,1,45250,0),
  @(45648,0,45642,0),
  @(45655,1,45649,0)
)

$dateFmt = $ws1.Range("A2").NumberFormat

$r = 2
foreach ($row in $sheet1Data) {
  $ws1.Cells.Item($r, 1).Value = $row[0]
  $ws1.Cells.Item($r, 1).NumberFormat = $dateFmt
  $ws1.Cells.Item($r, 2).Value = $row[1]
  $ws1.Cells.Item($r, 3).Value = $row[2]
  $ws1.Cells.Item($r, 3).NumberFormat = $dateFmt
  $ws1.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 2 - "Weekly Growth": header row only
# ---------------------------------------------------------------------------
$ws2.Range("A1").Value = "ds"
$ws2.Range("B1").Value = "PO_Requested_Qty"
$ws2.Range("C1").Value = "Growth%"

$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Sheet 3 - "Volume Insights": header row + a single summary row
# ---------------------------------------------------------------------------
$ws3.Range("A1").Value = "Total_PO_Quantity"
$ws3.Range("B1").Value = "Average_PO_Quantity"
$ws3.Range("C1").Value = "Max_PO_Quantity"
$ws3.Range("D1").Value = "Min_PO_Quantity"

$ws1.Range("A1:C1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

$ws3.Range("A2").Value = 0
$ws3.Range("B2").Value = 0
$ws3.Range("C2").Value = 0
$ws3.Range("D2").Value = 0

# ---------------------------------------------------------------------------
# Sheet 4 - "Prediction Info": header row + single predicted value
# ---------------------------------------------------------------------------
$ws4.Range("A1").Value = "Predicted_Next_Week_PO_Quantity"

$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)

$ws4.Range("A2").Value = 0

# ---------------------------------------------------------------------------
# Leave the first sheet active/selected, matching the original workbook
# ---------------------------------------------------------------------------
$ws1.Activate()
